$d = $word.ActiveDocument

# ------------------------------------------------------------------
# The amount in words needs to change from
#   "156.00(Trescientos Cuatro pesos 00/100 M.N.)"
# to
#   "156.00(Ciento Cincuenta y Seis 00/100 M.N.)"
# The "156" figure itself, and the trailing "00/100 M.N.)" stay the
# same; only the spelled-out amount changes, and the word "pesos" is
# dropped. There is a "_GoBack" bookmark sitting right after the
# figures that needs to end up positioned right after "...Cincuenta y S"
# (i.e. right before "eis") once the edit is done.
# ------------------------------------------------------------------

# 1) Replace the spelled-out amount. We search for the unique phrase
#    "Trescientos Cuatro pesos" (this purposely excludes the "156" run
#    that precedes it, so the existing "_GoBack" bookmark - anchored
#    right after "156" - is left untouched by this replacement).
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute(
    "Trescientos Cuatro pesos",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "Ciento Cincuenta y Seis",
    2
)

# 2) Relocate the "_GoBack" bookmark so it sits right after the new
#    "...Cincuenta y S" text (i.e. immediately before "eis"), matching
#    where the cursor was left by the edit.
$locate = $d.Content
$locate.Find.ClearFormatting()
$locate.Find.Execute(
    "Cincuenta y S",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "",
    0
)
$bookmarkPos = $locate.End
$d.Bookmarks.Add("_GoBack", $d.Range($bookmarkPos, $bookmarkPos))
